$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
# B1 already carried the bold/centered/bordered style (style index 1);
# give it its real label and extend the same formatting across the newly
# added header cells C1:F1.
$ws.Range("B1").Value = "现货变化"
$ws.Range("C1").Value = "代码"
$ws.Range("D1").Value = "期货价格"
$ws.Range("E1").Value = "180极限"
$ws.Range("F1").Value = "基差*极限"

$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1:F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data row (row 2) ----------------------------------------------------
# B2 (futures basis) already held its number; add the remaining columns.
# C2 ("1804") looks numeric but must stay a text code, so build it as a
# text formula result and paste the *value* in, which keeps it text
# without forcing a quote-prefixed number format on the cell.
$ws.Range("Z1").Formula = "=""1804"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4163) | Out-Null     # xlPasteValues
$ws.Range("Z1").ClearContents()
$excel.CutCopyMode = 0

$ws.Range("D2").Value = 14385
$ws.Range("E2").Value = 0.3868613138686132
$ws.Range("F2").Value = -0.7260499162548335
